# Added review for the coding phase
$wb = $excel.ActiveWorkbook

$reqWs  = $wb.Worksheets.Item("Requirements Phase Defects")
$archWs = $wb.Worksheets.Item("Architect. Design Phase Defects")
$codeWs = $wb.Worksheets.Item("Coding Phase Defects")

# --- Fill in the reviewer block (I3:J4) on the Coding Phase Defects sheet ---
$codeWs.Range("I3").Value = "Marginean Tiberius"
$codeWs.Range("J3").Value = 234
$codeWs.Range("I4").Value = "Mag Raul"
$codeWs.Range("J4").Value = 234

# --- Fill in the review table rows 10-14 (Crt.No / Checked Item / Doc. page-line / Comments) ---
$codeWs.Range("C10").Value = "C04"
$codeWs.Range("D10").Value = "OrdersGUIController / 88, 92`nKitchenGUIController / 55, 64"
$codeWs.Range("E10").Value = "Formatul de afisare al orei este gresit in cazul in care minutele sunt mai mici de 10 (ex: 13:05 va fi afisat 13:5)"

$codeWs.Range("C11").Value = "C06"
$codeWs.Range("D11").Value = "OrdersGUIController / 141`nKitchenGUIContorller / 53, 61"
$codeWs.Range("E11").Value = "Aplicatia permite apasarea unor butoane fara a selecta un element dintr-o lista in situatii in care ar fi necesar acest lucru. Astfel apar erori netratate"

$codeWs.Range("C12").Value = "C01"
$codeWs.Range("D12").Value = "MainGUIController`nKitchenGUIContorller"
$codeWs.Range("E12").Value = "Nu se asteapta inchiderea tuturor meselor inainte de inchiderea bucatariei, respectiv inchiderea bucatariei inainte de inchiderea restaurantului"

$codeWs.Range("C13").Value = "C09"
$codeWs.Range("D13").Value = "MenuDataModel / 12"
$codeWs.Range("E13").Value = "Numele parametrilor constructorului creeaza confuzie"

$codeWs.Range("C14").Value = "C03"
$codeWs.Range("D14").Value = "KitchenGUIController / 29"
$codeWs.Range("E14").Value = "Bucla while continua executia pana la inchiderea aplicatiei, chiar daca fereastra aferenta bucatariei a fost inchisa mai devreme"

# Row heights grow to fit the wrapped, multi-line comments that were just entered
$codeWs.Rows.Item(10).RowHeight = 30
$codeWs.Rows.Item(11).RowHeight = 45
$codeWs.Rows.Item(12).RowHeight = 45
$codeWs.Rows.Item(14).RowHeight = 45

# Effort to review document (hours)
$codeWs.Range("E32").Value = 1.5

# Widen columns D, E and H on the Coding Phase Defects sheet to fit the new text
$codeWs.Columns.Item(4).ColumnWidth = 27.666666666666668
$codeWs.Columns.Item(5).ColumnWidth = 54.333333333333336
$codeWs.Columns.Item(8).ColumnWidth = 10.5

# --- Update selections / active sheet ---
# Requirements sheet selection is untouched; Architect. Design sheet loses the
# tab-selected flag but keeps an updated selected cell; Coding Phase Defects
# becomes the newly active/selected tab.
$archWs.Range("I11").Select() | Out-Null

$codeWs.Activate()
$codeWs.Range("I31").Select() | Out-Null
